# Practice Task N8 complete
# Adds a new data row (A9/B9) mirroring the existing "Грабли" entry with a
# new price, then leaves the active selection on the newly entered cell -
# matching the manual data-entry workflow the author performed in Excel/Calc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Грабли"
$ws.Range("B9").Value = 786

$ws.Range("B9").Select()
